$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the misspelled name "Leeory" -> "Leeroy".
#    In the real edit the correction was made in-place inside Word, which
#    split the word into several runs ("Lee" | "r" | "o" | "y") and left the
#    "_GoBack" bookmark (marking the last edit position) sitting right
#    before the final "y". We reproduce both the corrected text and that
#    exact run/bookmark layout.
# ---------------------------------------------------------------------------

$nameRange = $d.Content
$nameRange.Find.Execute("Leeory", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

$start = $nameRange.Start
$end = $nameRange.End

# Replace the whole word first (same length: 6 chars -> 6 chars, so nothing
# downstream shifts).
$wordRange = $d.Range($start, $end)
$wordRange.Text = "Leeroy"

# Force run boundaries after "Lee" and after "Leero" by dropping (and
# immediately removing) throw-away bookmarks at those two offsets; a
# bookmark start/end can only live between runs, so OOXML has to split the
# text run around it. The split persists even after the temporary bookmark
# is deleted.
$d.Bookmarks.Add("zLeeroySplit1", $d.Range($start + 3, $start + 3)) | Out-Null
$d.Bookmarks.Add("zLeeroySplit2", $d.Range($start + 4, $start + 4)) | Out-Null

# The real "_GoBack" bookmark belongs right before the final "y" (offset
# start+5). Word only ever keeps a single "_GoBack" bookmark, so (re)adding
# it under that name automatically relocates it here and removes the one
# that used to sit at the very end of the document.
$d.Bookmarks.Add("_GoBack", $d.Range($start + 5, $start + 5)) | Out-Null

$d.Bookmarks.Item("zLeeroySplit1").Delete()
$d.Bookmarks.Item("zLeeroySplit2").Delete()

# ---------------------------------------------------------------------------
# 2) Merge "persiguiendo nuevas aventuras" + "!" + " " into a single run,
#    dropping the grammar-check markers that used to wrap the "!".
# ---------------------------------------------------------------------------

$phraseRange = $d.Content
$phraseRange.Find.Execute("persiguiendo nuevas aventuras! ", $true, $false, `
                           $false, $false, $false, $true, 1, $false, "", 0)

$pStart = $phraseRange.Start
$pEnd = $phraseRange.End

# Pin the two boundaries with throw-away bookmarks so the upcoming text
# replace only coalesces runs *inside* the span, instead of bleeding into
# the neighbouring runs (the preceding " " and the following "Tu tarea...").
$d.Bookmarks.Add("zPhraseStart", $d.Range($pStart, $pStart)) | Out-Null
$d.Bookmarks.Add("zPhraseEnd", $d.Range($pEnd, $pEnd)) | Out-Null

# The final text is already what those three runs spell out together, so a
# direct re-assignment would be a content no-op and the runs would be left
# untouched. Round-trip through a placeholder of identical length to force
# the engine to actually rebuild (and thus coalesce) the run(s) in range.
$placeholder = ""
for ($i = $pStart; $i -lt $pEnd; $i++) { $placeholder += "*" }
$phraseRange2 = $d.Range($pStart, $pEnd)
$phraseRange2.Text = $placeholder

$phraseRange3 = $d.Range($pStart, $pEnd)
$phraseRange3.Text = "persiguiendo nuevas aventuras! "

$d.Bookmarks.Item("zPhraseStart").Delete()
$d.Bookmarks.Item("zPhraseEnd").Delete()
